# Generate Report for Handoff
# - Overview sheet, rows 4-7 ("Latest HO Xliff Generate Date", col G) and
#   de-de sheet, rows 4-7 ("Latest Handoff Datetime", col H) share the same
#   underlying string and both advance from 2016-08-29 16:34:15 -> 16:34:40
# - zh-cn / de-de sheets: rows 4-7 "Priority" column (E) moves from "low" to "ht"
# - zh-cn sheet: rows 4-7 "Latest Handoff Datetime" column (H) advances from
#   2016-08-29 16:34:10 -> 2016-08-29 16:34:35

$wb = $excel.ActiveWorkbook

$overview = $wb.Sheets.Item("Overview")
$overview.Range("G4").Value = "2016-08-29 16:34:40"
$overview.Range("G5").Value = "2016-08-29 16:34:40"
$overview.Range("G6").Value = "2016-08-29 16:34:40"
$overview.Range("G7").Value = "2016-08-29 16:34:40"

$zhcn = $wb.Sheets.Item("zh-cn")
$zhcn.Range("E4").Value = "ht"
$zhcn.Range("E5").Value = "ht"
$zhcn.Range("E6").Value = "ht"
$zhcn.Range("E7").Value = "ht"
$zhcn.Range("H4").Value = "2016-08-29 16:34:35"
$zhcn.Range("H5").Value = "2016-08-29 16:34:35"
$zhcn.Range("H6").Value = "2016-08-29 16:34:35"
$zhcn.Range("H7").Value = "2016-08-29 16:34:35"

$dede = $wb.Sheets.Item("de-de")
$dede.Range("E4").Value = "ht"
$dede.Range("E5").Value = "ht"
$dede.Range("E6").Value = "ht"
$dede.Range("E7").Value = "ht"
$dede.Range("H4").Value = "2016-08-29 16:34:40"
$dede.Range("H5").Value = "2016-08-29 16:34:40"
$dede.Range("H6").Value = "2016-08-29 16:34:40"
$dede.Range("H7").Value = "2016-08-29 16:34:40"
